# Append a new "2025-03-09" price row (row 8) to every price sheet in the
# workbook, carrying forward the same value that was already in row 7 of
# that sheet (matches the previous day's price being repeated).
#
# Values (row 8, column B) keyed by sheet name, taken from row 7 of each
# sheet:
#   N-Dense                    -> 40
#   N-Type                     -> 43
#   N-type Wafer                -> 1.19
#   Cell Topcon 183mm           -> 0.295
#   Module Topcon 183mm         -> 0.1
#   Silver Rear_side            -> 5,263
#   Silver Busbar front-side    -> 7,879
#   Silver finger front-side    -> 7,929
#   USD_CNY                     -> 7.2647

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-09"

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.295"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,263"
    "Silver Busbar front-side"   = "7,879"
    "Silver finger front-side"   = "7,929"
    "USD_CNY"                    = "7.2647"
}

foreach ($ws in $wb.Worksheets) {
    $price = $sheetValues[$ws.Name]
    if ($price -eq $null) {
        continue
    }

    # Find the last used row in column A (the existing data ends at row 7
    # for every sheet before this edit), then append directly below it.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    $dateCell = $ws.Cells.Item($newRow, 1)
    $priceCell = $ws.Cells.Item($newRow, 2)

    # Force text storage (matching the rest of the column, which is stored
    # as text rather than as real numbers/dates) by marking the cell as
    # Text before assigning the value, then reset the style back to
    # Normal so no stray number-format style lingers on the cell.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate
    $dateCell.Style = "Normal"

    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"
}
